# Update the "想去人数" (want-to-attend count) values in column F
# for both the "展览" and "全部类型" worksheets.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 6673
    3  = 191
    5  = 50
    6  = 2075
    7  = 1580
    10 = 463
    12 = 5650
    13 = 77
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
